# SAO_SEPE.xlsx automatic update:
#  - rename "Paineis DARQ" -> "PAINEIS DARQ"
#  - rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#  - remove the "Desarquivamentos Pendentes" sheet

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Paineis DARQ")
$ws1.Name = "PAINEIS DARQ"

$ws6 = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$ws6.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# suppress the "permanently delete this sheet" confirmation dialog
$excel.DisplayAlerts = $false
$ws8 = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$ws8.Delete()
$excel.DisplayAlerts = $true

# keep the dashboard sheet as the active tab, as it was originally
$ws1.Activate()
